$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SE Palmeiras vs EC Vitoria Salvador)
$ws.Range("F2").Value = 1.32
$ws.Range("H2").Value = 11
$ws.Range("J2").Value = 5.6
$ws.Range("W2").Value = 3.8

# Row 4 (Atletico Bucaramanga vs Santa Fe)
$ws.Range("G4").Value = 2.14
$ws.Range("H4").Value = 3.9
$ws.Range("J4").Value = 3.05
$ws.Range("Q4").Value = 2.26
